$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.7
$ws.Range("I2").Value = 2.15
$ws.Range("K2").Value = 1.95
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("W2").Value = 8.5
$ws.Range("AH2").Value = 9
$ws.Range("AI2").Value = 9.5
$ws.Range("AN2").Value = 5.5
$ws.Range("BD2").Value = 126

# Row 3
$ws.Range("G3").Value = 1.9
$ws.Range("I3").Value = 4.75
$ws.Range("L3").Value = 5.5
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.53
$ws.Range("W3").Value = 5
$ws.Range("AX3").Value = 29

# Row 5
$ws.Range("G5").Value = 8.5
$ws.Range("I5").Value = 1.45
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 8

# Row 6
$ws.Range("O6").Value = 1.25
$ws.Range("P6").Value = 3.75
$ws.Range("Q6").Value = 1.8
$ws.Range("R6").Value = 2
